# "update opnam post data"
# The sample data row for "TEST123 / Testing Opnam" (row 2) is removed, and
# the remaining data row ("01023-50622 / BOLT, SEMS", previously row 3) is
# moved up to row 2 and updated with a new part number/name
# ("967120Z000 / HINGE LH"), keeping its existing Qty/Harga (5 / 15000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 2 ("TEST123" / "Testing Opnam") - this shifts the former
# row 3 ("01023-50622" / "BOLT, SEMS", Pcs, 5, 15000) up into row 2,
# carrying its row height and cell styles along with it.
$ws.Rows(2).Delete()

# Row 2 "No" column should restart the numbering at 1.
$ws.Range("A2").Value = 1

# New Part Number for the remaining row. The old value ("01023-50622") was
# stored with a quote-prefix (text) style; the new part number does not need
# that, so reset the cell back to the workbook's default style.
$ws.Range("B2").Value = "967120Z000"
$ws.Range("B2").Style = "Normal"

# New Part Name for the remaining row - keep its existing style.
$ws.Range("C2").Value = "HINGE LH"

# Uom (D2), Actual Stock (E2=5) and Harga Satuan (F2=15000) already carry the
# correct values after the row shift above, so nothing else to change there.

# Page setup now explicitly records a portrait orientation.
$ws.PageSetup.Orientation = 1

# Keep the active selection in sync with the new last cell (old F3, now F2).
$null = $ws.Range("F2").Select()
